# Fruta / hortaliza, semanal
#
# A new daily price record (Crimpson Seedless / Segunda) is inserted just
# above the existing "Autumn Royal" record for Macroferia Regional de Talca
# - Uva, pushing every following record down by one row
# (old row 704 -> new row 705, ... old row 795 -> new row 796).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 704; Excel shifts rows 704:795 down to 705:796.
$ws.Rows(704).Insert()

# Populate the freshly inserted row 704 with the new record. The
# "dimension"-only columns (A,B,C,E,F,G,H,I,J,R) repeat the same market /
# product metadata found on every other row of this sheet.
$ws.Range("A704").Value2 = 5
$ws.Range("B704").Value2 = "Macroferia Regional de Talca"
$ws.Range("C704").Value2 = "Maule"
$ws.Range("D704").Value2 = 45154
$ws.Range("E704").Value2 = 7
$ws.Range("F704").Value2 = "Fruta"
$ws.Range("G704").Value2 = 100109
$ws.Range("H704").Value2 = "Uva"
$ws.Range("I704").Value2 = 100109001
$ws.Range("J704").Value2 = "Uva"
$ws.Range("K704").Value2 = "Crimpson Seedless"
$ws.Range("L704").Value2 = "Segunda"
$ws.Range("M704").Value2 = 600
$ws.Range("N704").Value2 = 12000
$ws.Range("O704").Value2 = 12000
$ws.Range("P704").Value2 = 12000
$ws.Range("Q704").Value2 = "`$/bandeja 8 kilos"
$ws.Range("R704").Value2 = "Región de O'Higgins"
$ws.Range("S704").Value2 = 1500
$ws.Range("T704").Value2 = 8
